# Update the "Password recover" / reset-password translation rows on the
# "Worksheet" sheet.  Row 11 held placeholder text "Pug" in the English,
# German and French columns; row 12 held the same placeholder.  They are
# replaced with the real translations (English/French stay the same as the
# already-correct Portuguese column, German gets a proper translation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Row 11: "Password recover"
$ws.Range("B11").Value = "Password recover"
$ws.Range("C11").Value = "Passwort-Wiederherstellung"
$ws.Range("D11").Value = "Password recover"
$ws.Range("F11").Value = "Password recover"

# Row 12: "Your password reseted successully!"
$ws.Range("B12").Value = "Your password reseted successully!"
$ws.Range("C12").Value = "Ihr Passwort wurde erfolgreich zurückgesetzt!"
$ws.Range("D12").Value = "Your password reseted successully!"
$ws.Range("F12").Value = "Your password reseted successully!"
